# Updates cryptos list (prices + 1h volume %) per latest scrape; also
# reflects the two rank swaps (TheGraph/Dai at rows 38-39, dogwifhat/Maker
# at rows 43-44) produced by the upstream ranking refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Column D holds price strings that often parse as plain numbers
    # (e.g. "575.00", "0.120"). Assigning through .Value on a General-
    # formatted cell lets Excel coerce them to numeric, dropping the
    # trailing zeros / exact text the source data relies on. Stash the
    # cell's existing style, force Text format for the write, then put
    # the original style back so no formatting actually changes.
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") '69.529.93'
$ws.Range("E2").Value = '  -0.75%  '
Set-TextValue $ws.Range("D3") '3.506.30'
$ws.Range("E3").Value = '  -2.34%  '
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.14%  '
Set-TextValue $ws.Range("D5") '575.00'
$ws.Range("E5").Value = '  -0.88%  '
Set-TextValue $ws.Range("D6") '186.42'
$ws.Range("E6").Value = '  -2.40%  '
Set-TextValue $ws.Range("D7") '3.498.83'
$ws.Range("E7").Value = '  -2.46%  '
Set-TextValue $ws.Range("D8") '0.612'
$ws.Range("E8").Value = '  -3.25%  '
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("E10").Value = '  +4.01%  '
Set-TextValue $ws.Range("D11") '0.648'
$ws.Range("E11").Value = '  -2.49%  '
Set-TextValue $ws.Range("D12") '54.15'
$ws.Range("E12").Value = '  -3.08%  '
Set-TextValue $ws.Range("D13") '0.0000302'
$ws.Range("E13").Value = '  -1.74%  '
Set-TextValue $ws.Range("D14") '9.45'
$ws.Range("E14").Value = '  -2.52%  '
Set-TextValue $ws.Range("D15") '4.054.02'
$ws.Range("E15").Value = '  -2.79%  '
Set-TextValue $ws.Range("D17") '69.345.41'
$ws.Range("E17").Value = '  -1.02%  '
Set-TextValue $ws.Range("D18") '3.495.26'
$ws.Range("E18").Value = '  -2.61%  '
Set-TextValue $ws.Range("D19") '12.25'
$ws.Range("E19").Value = '  -3.66%  '
Set-TextValue $ws.Range("D20") '0.120'
$ws.Range("E20").Value = '  -1.05%  '
Set-TextValue $ws.Range("D21") '540.60'
$ws.Range("E21").Value = '  +12.49%  '
$ws.Range("E22").Value = '  -3.57%  '
Set-TextValue $ws.Range("D23") '18.41'
$ws.Range("E23").Value = '  -3.63%  '
Set-TextValue $ws.Range("D24") '4.97'
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("E25").Value = '  +0.55%  '
Set-TextValue $ws.Range("D26") '93.62'
$ws.Range("E26").Value = '  -2.27%  '
Set-TextValue $ws.Range("D27") '11.29'
$ws.Range("E27").Value = '  +1.51%  '
Set-TextValue $ws.Range("D28") '2.95'
$ws.Range("E28").Value = '  -1.74%  '
Set-TextValue $ws.Range("D29") '9.13'
$ws.Range("E29").Value = '  -3.28%  '
Set-TextValue $ws.Range("D30") '31.83'
$ws.Range("E30").Value = '  -0.99%  '
Set-TextValue $ws.Range("D31") '7.26'
$ws.Range("E31").Value = '  -5.19%  '
Set-TextValue $ws.Range("D32") '12.64'
$ws.Range("E32").Value = '  +3.37%  '
Set-TextValue $ws.Range("D33") '64.53'
$ws.Range("E33").Value = '  -3.26%  '
Set-TextValue $ws.Range("D34") '0.114'
$ws.Range("E34").Value = '  -5.31%  '
Set-TextValue $ws.Range("D35") '535.51'
$ws.Range("E35").Value = '  -8.33%  '
Set-TextValue $ws.Range("D36") '3.10'
$ws.Range("E36").Value = '  +9.08%  '
Set-TextValue $ws.Range("D37") '37.98'
$ws.Range("E37").Value = '  -2.80%  '
$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D38") '0.401'
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D39") '1.00'
$ws.Range("E39").Value = '  -0.06%  '
Set-TextValue $ws.Range("D40") '0.0₃0764'
$ws.Range("E40").Value = '  -5.02%  '
Set-TextValue $ws.Range("D41") '3.38'
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("E42").Value = '  -2.37%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D43") '3.05'
$ws.Range("E43").Value = '  -6.98%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D44") '3.293.64'
$ws.Range("E44").Value = '  +2.35%  '
Set-TextValue $ws.Range("D45") '2.98'
$ws.Range("E45").Value = '  -2.89%  '
Set-TextValue $ws.Range("D46") '0.0444'
$ws.Range("E46").Value = '  -1.38%  '
Set-TextValue $ws.Range("D47") '3.48'
$ws.Range("E47").Value = '  +4.31%  '
Set-TextValue $ws.Range("D48") '0.134'
$ws.Range("E48").Value = '  -3.28%  '
Set-TextValue $ws.Range("D49") '8.91'
$ws.Range("E49").Value = '  -6.47%  '
Set-TextValue $ws.Range("D50") '0.997'
$ws.Range("E50").Value = '  -0.16%  '
Set-TextValue $ws.Range("D51") '137.18'
$ws.Range("E51").Value = '  +1.63%  '
